$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cardholder name
$ws.Range("C2").Value = "Hartmut"

# Card number (kept as text) and surname
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats: reapply text formatting (no quote-prefix residue)
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 16.05.2024"

$ws.Range("B6").Value = "19.05."
$ws.Range("C6").Value = "20.05."
$ws.Range("D6").Value = "PAYPAL LZTSUB"
$ws.Range("E6").Value = "33,00-"

$ws.Range("B7").Value = "20.05."
$ws.Range("C7").Value = "21.05."
$ws.Range("D7").Value = "PAYPAL SUAQAF"
$ws.Range("E7").Value = "95,37-"

$ws.Range("B8").Value = "21.05."
$ws.Range("C8").Value = "22.05."
$ws.Range("D8").Value = "KARTENZ./21.05 REWE RO"
$ws.Range("E8").Value = "129,25-"

$ws.Range("B9").Value = "25.05."
$ws.Range("C9").Value = "26.05."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 89330599"
$ws.Range("E9").Value = "40,15-"

$ws.Range("B10").Value = "27.05."
$ws.Range("C10").Value = "28.05."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-77143014"
$ws.Range("E10").Value = "56,21-"

$ws.Range("B11").Value = "30.05."
$ws.Range("C11").Value = "31.05."
$ws.Range("D11").Value = "PAYPAL YFMMBL"
# E11 was a previously-empty cell with a different style (s=12); the new row
# matches the formatting used by the other amount cells (E6:E10, s=17).
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E11").Value = "65,04-"

# Closing balance line and next billing date
$ws.Range("D12").Value = "KONTOSTAND AM 01.06.2024"
$ws.Range("E12").Value = "419,02-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 08.06.2024"
